$d = $word.ActiveDocument

$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Locate the "Axios" and "Formik" paragraphs by content rather than a fixed
# index, so the script is resilient to the paragraph collection positions.
$axiosIndex = 0
$formikIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text.Trim()
    if ($t -eq "Axios") { $axiosIndex = $i }
    if ($t -eq "Formik") { $formikIndex = $i }
}

# --- "Axios" -> "Axios (https://github.com/axios/axios)" ---
# Drop the paragraph's rtl pPr, keep the spell-checked "Axios" run, append
# " (", the URL, ")" as separate runs, and move the _GoBack bookmark to the
# end of this paragraph.
$axiosXml = '<w:p ' + $wns + '>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>Axios</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> (</w:t></w:r>' +
  '<w:r><w:t>https://github.com/axios/axios</w:t></w:r>' +
  '<w:r><w:t>)</w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
  '<w:bookmarkEnd w:id="0"/>' +
  '</w:p>'

$axiosPara = $d.Paragraphs.Item($axiosIndex)
$axiosPara.Range.InsertXML($axiosXml)

# --- "Formik" -> "Formik (https://github.com/jaredpalmer/formik)" ---
# No bookmark anymore (it moved to the Axios paragraph above); add spell
# proofing marks around "Formik" and append the " (url)" runs.
$formikXml = '<w:p ' + $wns + '>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>Formik</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> (</w:t></w:r>' +
  '<w:r><w:t>https://github.com/jaredpalmer/formik</w:t></w:r>' +
  '<w:r><w:t>)</w:t></w:r>' +
  '</w:p>'

$formikPara = $d.Paragraphs.Item($formikIndex)
$formikPara.Range.InsertXML($formikXml)

# --- New paragraph after "Formik": "React-Toastify (https://github.com/fkhadra/react-toastify)" ---
$formikPara = $d.Paragraphs.Item($formikIndex)
$formikPara.Range.InsertParagraphAfter()

$toastifyXml = '<w:p ' + $wns + '>' +
  '<w:r><w:t>React-</w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>Toastify</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> (</w:t></w:r>' +
  '<w:r><w:t>https://github.com/fkhadra/react-toastify</w:t></w:r>' +
  '<w:r><w:t>)</w:t></w:r>' +
  '</w:p>'

$toastifyPara = $d.Paragraphs.Item($formikIndex + 1)
$toastifyPara.Range.InsertXML($toastifyXml)
